$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '89.637.19'
$ws.Cells.Item(2, 5).Value = '  -0.92%  '

$ws.Cells.Item(3, 4).Value = '3.066.29'
$ws.Cells.Item(3, 5).Value = '  -2.20%  '

$ws.Cells.Item(4, 5).Value = '  -0.17%  '

$ws.Cells.Item(5, 4).Value = '236.33'
$ws.Cells.Item(5, 5).Value = '  +9.79%  '

$ws.Cells.Item(6, 4).Value = '616.59'
$ws.Cells.Item(6, 5).Value = '  -0.75%  '

$ws.Cells.Item(7, 4).Value = '1.07'
$ws.Cells.Item(7, 5).Value = '  -5.93%  '

$ws.Cells.Item(8, 5).Value = '  +0.06%  '

$ws.Cells.Item(9, 5).Value = '  -0.04%  '

$ws.Cells.Item(10, 4).Value = '3.065.97'
$ws.Cells.Item(10, 5).Value = '  -2.16%  '

$ws.Cells.Item(11, 4).Value = '0.706'
$ws.Cells.Item(11, 5).Value = '  -5.17%  '

$ws.Cells.Item(12, 4).Value = '0.199'
$ws.Cells.Item(12, 5).Value = '  -0.53%  '

$ws.Cells.Item(13, 5).Value = '  +0.32%  '

$ws.Cells.Item(14, 4).Value = '34.63'
$ws.Cells.Item(14, 5).Value = '  -1.15%  '

$ws.Cells.Item(15, 4).Value = '89.303.67'
$ws.Cells.Item(15, 5).Value = '  -1.06%  '

$ws.Cells.Item(16, 4).Value = '5.35'
$ws.Cells.Item(16, 5).Value = '  -5.29%  '

$ws.Cells.Item(17, 4).Value = '3.637.45'
$ws.Cells.Item(17, 5).Value = '  -2.00%  '

$ws.Cells.Item(18, 4).Value = '3.061.38'
$ws.Cells.Item(18, 5).Value = '  -2.54%  '

$ws.Cells.Item(19, 4).Value = '3.75'
$ws.Cells.Item(19, 5).Value = '  +1.21%  '

$ws.Cells.Item(20, 4).Value = '0.0000212'
$ws.Cells.Item(20, 5).Value = '  -0.44%  '

$ws.Cells.Item(21, 4).Value = '13.74'
$ws.Cells.Item(21, 5).Value = '  -5.15%  '

$ws.Cells.Item(22, 4).Value = '429.53'
$ws.Cells.Item(22, 5).Value = '  -7.15%  '

$ws.Cells.Item(23, 4).Value = '5.39'
$ws.Cells.Item(23, 5).Value = '  +1.99%  '

$ws.Cells.Item(24, 4).Value = '8.66'
$ws.Cells.Item(24, 5).Value = '  -4.20%  '

$ws.Cells.Item(25, 4).Value = '5.55'
$ws.Cells.Item(25, 5).Value = '  -6.13%  '

$ws.Cells.Item(26, 4).Value = '86.73'
$ws.Cells.Item(26, 5).Value = '  -8.58%  '

$ws.Cells.Item(27, 4).Value = '11.64'
$ws.Cells.Item(27, 5).Value = '  -4.82%  '

$ws.Cells.Item(28, 4).Value = '3.232.27'
$ws.Cells.Item(28, 5).Value = '  -2.45%  '

$ws.Cells.Item(29, 5).Value = '  +0.17%  '

$ws.Cells.Item(30, 5).Value = '  +41.29%  '

$ws.Cells.Item(31, 2).Value = 'Cronos'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(31, 4).Value = '0.157'
$ws.Cells.Item(31, 5).Value = '  -4.03%  '

$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).Value = '8.93'
$ws.Cells.Item(32, 5).Value = '  -2.64%  '

$ws.Cells.Item(33, 4).Value = '0.199'
$ws.Cells.Item(33, 5).Value = '  -5.78%  '

$ws.Cells.Item(34, 4).Value = '25.53'
$ws.Cells.Item(34, 5).Value = '  -4.35%  '

$ws.Cells.Item(35, 5).Value = '  +2.26%  '

$ws.Cells.Item(36, 4).Value = '3.99'
$ws.Cells.Item(36, 5).Value = '  +56.18%  '

$ws.Cells.Item(37, 4).Value = '7.11'
$ws.Cells.Item(37, 5).Value = '  +1.91%  '

$ws.Cells.Item(38, 4).Value = '490.43'
$ws.Cells.Item(38, 5).Value = '  -5.01%  '

$ws.Cells.Item(39, 4).Value = '3.61'
$ws.Cells.Item(39, 5).Value = '  +0.59%  '

$ws.Cells.Item(40, 5).Value = '  -3.09%  '

$ws.Cells.Item(41, 4).Value = '0.0902'
$ws.Cells.Item(41, 5).Value = '  -0.85%  '

$ws.Cells.Item(42, 4).Value = '1.25'
$ws.Cells.Item(42, 5).Value = '  -5.45%  '

$ws.Cells.Item(43, 4).Value = '22.09'
$ws.Cells.Item(43, 5).Value = '  -0.58%  '

$ws.Cells.Item(44, 5).Value = '  -0.01%  '

$ws.Cells.Item(45, 4).Value = '0.397'
$ws.Cells.Item(45, 5).Value = '  -6.41%  '

$ws.Cells.Item(46, 4).Value = '156.44'
$ws.Cells.Item(46, 5).Value = '  +4.11%  '

$ws.Cells.Item(47, 4).Value = '1.85'
$ws.Cells.Item(47, 5).Value = '  -6.42%  '

$ws.Cells.Item(48, 4).Value = '0.671'
$ws.Cells.Item(48, 5).Value = '  -7.87%  '

$ws.Cells.Item(49, 4).Value = '44.34'
$ws.Cells.Item(49, 5).Value = '  -2.25%  '

$ws.Cells.Item(50, 5).Value = '  -0.22%  '

$ws.Cells.Item(51, 5).Value = '  -4.99%  '

